$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Number formats already used by the rest of the table (same codes as the
# existing "numFmtId=180" date style and "numFmtId=178" amount style), so
# that new cells line up with style ids 2 and 3 used throughout the sheet.
$dateFormat = "[`$-409]d/mmm/yyyy;@"
$amountFormat = "`"₹`"#,##0;`"₹`"\-#,##0"

# New entries to append after the existing last row (331).
# Columns: A=DATE, B=VEHICLE REG NO, C=VEHICLE BRAND, D=ISSUE, E=STATUS, F=AMOUNT, G=CASH TYPE
# Date serials: 44802 = 29-Aug-2022, 44803 = 30-Aug-2022
$rows = @(
    @{ Row=332; Date=44802; B="KA03MQ5513"; C="HONDA CITY";  D="PMS";                    E="WORK DONE DELIVERED"; F=5558;  G="G PAY" },
    @{ Row=333; Date=44802; B="KA03MM7229"; C="I10";         D="GEAR SWIFTING PROBLEM";  E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=334; Date=44802; B="KA04MN7120"; C="FIESTA";      D="PMS";                    E="WORK DONE DELIVERED"; F=8095;  G="G PAY" },
    @{ Row=335; Date=44802; B="KA41P3439";  C="NEW VERNA";   D="RUNNING REPAIR";         E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=336; Date=44802; B="KA01ML2754"; C="ECOSPORT";    D="RUNNING REPAIR";         E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=337; Date=44802; B="KA03MU6336"; C="HONDA CITY";  D="BRAKE PAD CHANGE";       E="WORK DONE DELIVERED"; F=2000;  G="G PAY" },
    @{ Row=338; Date=44802; B="KA03MS1179"; C="DZIER";       D="RUNNING REPAIR";         E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=339; Date=44802; B="KA53N7601";  C="BOLERO";      D="FUEL PIPE CHANGE";       E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=340; Date=44803; B="KA05MJ6661"; C="I20";         D="BODY SHOP";              E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=341; Date=44803; B="KA53Z9119";  C="I10";         D="FENDER LINING CHANGE";   E="WORK DONE DELIVERED"; F=1721;  G=$null },
    @{ Row=342; Date=44803; B="KA02MA7199"; C="SCORPIO";     D="WIRING PROBLEM";         E="WORK IN PROGRESS";    F=$null; G=$null },
    @{ Row=343; Date=44803; B="KA03ND8314"; C="SWIFT";       D="PMS";                    E="WORK DONE DELIVERED"; F=12158; G=$null },
    @{ Row=344; Date=44803; B="KA03MT5176"; C="ALTO 800";    D="GENERAL CHECKUP";        E="WORK DONE DELIVERED"; F=1598;  G=$null },
    @{ Row=345; Date=44803; B="KA03MM7095"; C="HONDA JAZZ";  D="PMS";                    E="WORK DONE DELIVERED"; F=5500;  G="CASH" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E

    if ($null -ne $r.F) {
        $ws.Cells.Item($rowNum, 6).Value = $r.F
        $ws.Cells.Item($rowNum, 6).NumberFormat = $amountFormat
    }
    if ($null -ne $r.G) {
        $ws.Cells.Item($rowNum, 7).Value = $r.G
    }
}

# Reflect the final scroll position / selection from the author's edit.
$excel.ActiveWindow.ScrollRow = 322
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H345").Select()
